# #272 Ajout d'un scenario de recherche de l'offre d'un professionnel avec un ID Nat PS
#
# Two substantive changes in this StructureDefinition export:
#   1. The "Date" metadata value is bumped to the new publication timestamp.
#   2. On the "Elements" sheet, the two mapping columns - "Mapping: RIM Mapping"
#      (AK) and "Mapping: Spécification métier vers l'extension ROR
#      HealthcareServicePsychiatricSector" (AL) - trade places: AK becomes the
#      (wide) Spécification column and AL becomes the (narrow) RIM Mapping
#      column, for both the header row and every data row.

$wb = $excel.ActiveWorkbook

# 1. Bump the publication Date on the Metadata sheet (Property/Value table).
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# 2. Swap the AK ("Mapping: RIM Mapping") and AL ("Mapping: Spécification
#    métier...") columns on the Elements sheet: swap the header text and every
#    row's value between the two columns.
$wsEl = $wb.Worksheets.Item("Elements")

$lastRow = $wsEl.Range("A1").End(-4121).Row   # xlDown = -4121

for ($r = 1; $r -le $lastRow; $r++) {
    $akVal = $wsEl.Range("AK$r").Value2
    $alVal = $wsEl.Range("AL$r").Value2
    if ($akVal -ne $alVal) {
        $wsEl.Range("AK$r").Value = $alVal
        $wsEl.Range("AL$r").Value = $akVal
    }
}

# Swap the column widths to match: AK (now Spécification) becomes the wide
# column, AL (now RIM Mapping) becomes the narrow one.
$wsEl.Columns.Item(37).ColumnWidth = 91.15
$wsEl.Columns.Item(38).ColumnWidth = 24.16
